$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45 - MindSensorPressureSensor: claim the sensor (Developer/finished/Mode/interface)
$ws.Range("D45").Value = "Lawrie"
$ws.Range("E45").Value = "N"
$ws.Range("F45").Value = "Pressure"
$ws.Range("G45").Value = "SampleProvider"

# Row 69 - RCXTemperatureSensor: claim the sensor (Developer/finished/Mode/interface)
$ws.Range("D69").Value = "Lawrie"
$ws.Range("E69").Value = "N"
$ws.Range("F69").Value = "Temperature"
$ws.Range("G69").Value = "SampleProvider"

# Move the current selection (matches the saved view state in the workbook)
$ws.Range("I66").Select() | Out-Null
